# Assign to committee final changes , the navigations and saving
#
# Adds the new "Info" work item row to the pending work items list and
# moves the sheet selection, matching the author's final review pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pending work item (row 8): Type / Description / Date
$ws.Range("A8").Value = "Info"
$ws.Range("B8").Value = "Check migration EF , if its properly creating the planning meeting members and reviews relationship columns"
$ws.Range("C8").Value = 45971

# Final navigation before saving
$ws.Range("C15").Select()

$wb.Save()
